$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = "Record"
$ws.Cells.Item(4, 2).Value = "Balanço Geral"
$ws.Cells.Item(4, 3).Value = "Social"
$ws.Cells.Item(4, 4).Value = "2025-04-02T11:56"
$ws.Cells.Item(4, 5).Value = "Negativo"
$ws.Cells.Item(4, 6).Value = "Primeiro dia do mutirão de atualização do CadÚnico é marcado por tumulto. Repórter *ao vivo*. Grande quantidade de pessoas. 800 senhas distribuídas. Mutirão na Fundação de Esportes. Fila quase rodou o quarteirão. Pessoas chegaram de madrugada. Ninguém teria ido repassar informações para eles. Não teve prioridade de atendimento. Entre os entrevistados, mãe de autista e idosa com deficiência. Antes, atendimento era nos Cras, segundo entrevistada. Pessoas chamando de humilhação e covardia.*sem nota*"
